# 390-RBI-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment-Loanproduct4.xlsx
# "code refactoring and loan accounting and charges added"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # ProductLoan_Input
$ws2 = $wb.Worksheets.Item(2)   # ProductLoan_Output

# ---------------------------------------------------------------------------
# 1) Input sheet: a few existing values change
# ---------------------------------------------------------------------------

# shortname: was text "kar3" -> now numeric 390
$ws1.Range("B3").Value = 390

# nominalinterestratedefault: 12 -> 1
$ws1.Range("B11").Value = 1

# maximumallowedaoutstandingbalance: 5000 -> 10000
$ws1.Range("B28").Value = 10000

# ---------------------------------------------------------------------------
# 2) Append new loan-accounting rows 31-42 to the input sheet.
#    Shared-string order in the original commit put every *value* (column B)
#    before the corresponding *label* (column A), so we mirror that order:
#    write all of column B for the new rows first, then all of column A.
# ---------------------------------------------------------------------------

$ws1.Range("B31").Value = "Cash"
$ws1.Range("B32").Value = "Loan portfolio "
$ws1.Range("B33").Value = "Interest Receivable "
$ws1.Range("B34").Value = "Penalties Receivable "
$ws1.Range("B35").Value = "Transfer in Suspence "
$ws1.Range("B36").Value = "Fees Receivable"
$ws1.Range("B37").Value = "Income from interest"
$ws1.Range("B38").Value = "Income from penalties"
$ws1.Range("B39").Value = "Income from fees"
$ws1.Range("B40").Value = "Income from recovery repayments"
$ws1.Range("B41").Value = "Losses Writtenoff "
$ws1.Range("B42").Value = "Overpayment Liability"

$ws1.Range("A31").Value = "fundsource"
$ws1.Range("A32").Value = "loanprotfolio"
$ws1.Range("A33").Value = "interestreceivable"
$ws1.Range("A34").Value = "penaltiesreceivable"
$ws1.Range("A35").Value = "transferinsuspense"
$ws1.Range("A36").Value = "feesreceivable"
$ws1.Range("A37").Value = "incomefrominterest"
$ws1.Range("A38").Value = "incomefrompenalties"
$ws1.Range("A39").Value = "incomefromfees"
$ws1.Range("A40").Value = "incomefromrecoveryrepayments"
$ws1.Range("A41").Value = "loseswrittenoff"
$ws1.Range("A42").Value = "overpaymentliability"

# Re-apply the same formatting used by the rest of the sheet (column A uses
# the style of the "principaldefault"-style label cells, column B uses the
# style of the plain value cells) to the newly appended rows.
$ws1.Range("A9").Copy() | Out-Null
$ws1.Range("A31:A42").PasteSpecial(-4122) | Out-Null

$ws1.Range("B10").Copy() | Out-Null
$ws1.Range("B31:B42").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 3) Column B got a lot wider to fit the new descriptive text, and lost its
#    "best fit" auto-sizing flag (explicit width now).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 58.25

# ---------------------------------------------------------------------------
# 4) Selection / active-sheet bookkeeping: user scrolled the input sheet down
#    to see the new rows (selecting B30), then left the Output sheet active.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("B30").Select() | Out-Null

$ws2.Activate() | Out-Null
